$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 346.2
$ws.Range("I9").Value = 320.25
$ws.Range("K9").Value = 320.25
$ws.Range("M9").Value = -151.25

$ws.Range("H19").Value = 1102
$ws.Range("I19").Value = 1386
$ws.Range("K19").Value = 1386
$ws.Range("M19").Value = -1211

$ws.Range("H96").Value = 135.11111
$ws.Range("I96").Value = 93.5
$ws.Range("J96").Value = 218.33333
$ws.Range("K96").Value = 280.5
$ws.Range("L96").Value = 654.99999
$ws.Range("M96").Value = 1092.5
$ws.Range("N96").Value = -3400.99999

$ws.Range("H101").Value = 556
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = $null

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").Value = $null

$ws.Range("H107").Value = 318.75
$ws.Range("I107").Value = 157.41667
$ws.Range("K107").Value = 157.41667
$ws.Range("M107").Value = 1762.58333

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").Value = $null

$ws.Range("H135").Value = 531
$ws.Range("I135").Value = 582.5833
$ws.Range("K135").Value = 5243.2497
$ws.Range("M135").Value = -2708.2497

$ws.Range("H137").Value = 2684.9048
$ws.Range("I137").Value = 1992
$ws.Range("K137").Value = 5976
$ws.Range("M137").Value = -3426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 822.2
$ws.Range("I2").Value = 710.4
$ws.Range("K2").Value = 710.4
$ws.Range("M2").Value = -597.4

$ws.Range("H32").Value = 6256621.5
$ws.Range("I32").Value = 9769
$ws.Range("K32").Value = 9769
$ws.Range("M32").Value = -9482

$ws.Range("H45").Value = 3533.0908
$ws.Range("I45").Value = 2399.4
$ws.Range("J45").Value = 4477.8335
$ws.Range("K45").Value = 2399.4
$ws.Range("L45").Value = 4477.8335
$ws.Range("M45").Value = -2022.4
$ws.Range("N45").Value = -5231.8335

$ws.Range("H61").Value = 3286.4
$ws.Range("I61").Value = 3084.1667
$ws.Range("K61").Value = 3084.1667
$ws.Range("M61").Value = -2872.1667

$ws.Range("H74").Value = 1873.4722
$ws.Range("I74").Value = 1667.9259
$ws.Range("K74").Value = 1667.9259
$ws.Range("M74").Value = -793.9259

$ws.Range("H77").Value = 1873.4722
$ws.Range("I77").Value = 1667.9259
$ws.Range("K77").Value = 8339.629499999999
$ws.Range("M77").Value = -3971.629499999999

$ws.Range("H102").Value = 2755.5293
$ws.Range("I102").Value = 1789.6
$ws.Range("K102").Value = 1789.6
$ws.Range("M102").Value = -167.5999999999999

$ws.Range("H116").Value = 822.2
$ws.Range("I116").Value = 710.4
$ws.Range("K116").Value = 710.4
$ws.Range("M116").Value = 1583.6

$ws.Range("H122").Value = 1564
$ws.Range("I122").Value = 1138.5
$ws.Range("K122").Value = 3415.5
$ws.Range("M122").Value = -965.5

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = $null
$ws.Range("N128").Value = 0

$ws.Range("H136").Value = 3286.4
$ws.Range("I136").Value = 3084.1667
$ws.Range("K136").Value = 9252.500100000001
$ws.Range("M136").Value = -6702.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 822.2
$ws.Range("I3").Value = 710.4
$ws.Range("K3").Value = 710.4
$ws.Range("M3").Value = -596.4

$ws.Range("H80").Value = 660.0769
$ws.Range("I80").Value = 421.27274
$ws.Range("K80").Value = 421.27274
$ws.Range("M80").Value = 576.72726

$ws.Range("H83").Value = 660.0769
$ws.Range("I83").Value = 421.27274
$ws.Range("K83").Value = 2106.3637
$ws.Range("M83").Value = 2885.6363

$ws.Range("H105").Value = 2190.1428
$ws.Range("I105").Value = 2206.2
$ws.Range("K105").Value = 2206.2
$ws.Range("M105").Value = -459.1999999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6395.067
$ws.Range("I31").Value = 2467.125
$ws.Range("K31").Value = 2467.125
$ws.Range("M31").Value = -2172.125

$ws.Range("H34").Value = 6395.067
$ws.Range("I34").Value = 2467.125
$ws.Range("K34").Value = 2467.125
$ws.Range("M34").Value = -2265.125

$ws.Range("H105").Value = 1489.75
$ws.Range("I105").Value = 878.5714
$ws.Range("K105").Value = 878.5714
$ws.Range("M105").Value = 868.4286

$ws.Range("H132").Value = 3256.2856
$ws.Range("I132").Value = 2965.6667
$ws.Range("K132").Value = 8897.000100000001
$ws.Range("M132").Value = -6367.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 171924.25
$ws.Range("J37").Value = 171924.25
$ws.Range("L37").Value = 515772.75
$ws.Range("N37").Value = -515996.75

$ws.Range("H131").Value = 2367.7896
$ws.Range("J131").Value = 2543.1875
$ws.Range("L131").Value = 7629.5625
$ws.Range("N131").Value = -17709.5625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 78.55556
$ws.Range("I2").Value = 37.727272
$ws.Range("K2").Value = 37.727272
$ws.Range("M2").Value = 75.272728

$ws.Range("H59").Value = 28110
$ws.Range("J59").Value = 28110
$ws.Range("L59").Value = 28110
$ws.Range("N59").Value = -29276

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = $null
$ws.Range("N69").Value = 0

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = $null
$ws.Range("N72").Value = 0

$ws.Range("H80").Value = 5599.4
$ws.Range("I80").Value = 5599.4
$ws.Range("K80").Value = 5599.4
$ws.Range("M80").Value = -4601.4

$ws.Range("H83").Value = 5599.4
$ws.Range("I83").Value = 5599.4
$ws.Range("K83").Value = 27997
$ws.Range("M83").Value = -23005

$ws.Range("H86").Value = 100143
$ws.Range("J86").Value = 100143
$ws.Range("L86").Value = 100143
$ws.Range("N86").Value = -102515

$ws.Range("H89").Value = 100143
$ws.Range("J89").Value = 100143
$ws.Range("L89").Value = 300429
$ws.Range("N89").Value = -312285

$ws.Range("H102").Value = 1086.0588
$ws.Range("I102").Value = 1122.75
$ws.Range("K102").Value = 1122.75
$ws.Range("M102").Value = 499.25

$ws.Range("H114").Value = 125000
$ws.Range("J114").Value = 125000
$ws.Range("L114").Value = 125000
$ws.Range("N114").Value = -133678

$ws.Range("H126").Value = 3331.3333
$ws.Range("I126").Value = 2499.75
$ws.Range("J126").Value = 4994.5
$ws.Range("K126").Value = 7499.25
$ws.Range("L126").Value = 14983.5
$ws.Range("M126").Value = -5029.25
$ws.Range("N126").Value = -19923.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5000
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = $null

$ws.Range("H68").Value = 6888.778
$ws.Range("J68").Value = 10000
$ws.Range("L68").Value = 10000
$ws.Range("N68").Value = -11498

$ws.Range("H71").Value = 6888.778
$ws.Range("J71").Value = 10000
$ws.Range("L71").Value = 50000
$ws.Range("N71").Value = -57488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1399.3334
$ws.Range("I100").Value = 1227.8572
$ws.Range("K100").Value = 2455.7144
$ws.Range("M100").Value = -1914.7144

$ws.Range("H122").Value = 3808.8823
$ws.Range("I122").Value = 3796.9375
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 11390.8125
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -8940.8125
$ws.Range("N122").Value = -16900

$ws.Range("H126").Value = 6047.2104
$ws.Range("I126").Value = 3641.7144
$ws.Range("J126").Value = 7450.4165
$ws.Range("K126").Value = 10925.1432
$ws.Range("L126").Value = 22351.2495
$ws.Range("M126").Value = -8455.143199999999
$ws.Range("N126").Value = -27291.2495
